$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.499.99"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "'2.441.63"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.87%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'547.13"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "'145.75"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.582"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.39%  "
$ws.Range("D9").Value = "'2.439.48"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "'5.44"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.85%  "
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("D14").Value = "'25.83"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").Value = "'2.879.95"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("E16").Value = "  +4.04%  "
$ws.Range("D17").Value = "'61.644.57"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "'2.446.02"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("D19").Value = "'10.75"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.27%  "
$ws.Range("D20").Value = "'6.87"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("D22").Value = "'318.80"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  +6.53%  "
$ws.Range("D25").Value = "'63.73"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("D26").Value = "'0.0₃0970"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.76%  "
$ws.Range("D27").Value = "'2.567.13"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("D30").Value = "'7.79"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.41%  "
$ws.Range("D31").Value = "'525.41"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("D32").Value = "'8.18"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("E35").Value = "  +1.90%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  -2.25%  "
$ws.Range("D38").Value = "'4.72"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("E41").Value = "  +4.21%  "
$ws.Range("D42").Value = "'138.69"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.77%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "'40.33"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").Value = "'2.27"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("D46").Value = "'142.55"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.67%  "
$ws.Range("E47").Value = "  +1.85%  "
$ws.Range("D48").Value = "'21.17"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("D49").Value = "'0.0525"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("E51").Value = "  -0.87%  "
